# Edit script for backup_projetos.xlsx
# - Corrects a handful of "Feito"/"Nao Feito" status toggles for rows 106-131
# - Fills in previously-empty Mapa/Foto Perfil/Nome no Mapa columns (E/F/G) for rows 106-131
# - Appends 30 new task rows (380-409) for two new artists plus a small test project

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix Status Tarefa (column C) toggles for rows 106-131 ---
$ws.Cells.Item(110, 3).Value = 'Não Feito'
$ws.Cells.Item(111, 3).Value = 'Não Feito'
$ws.Cells.Item(112, 3).Value = 'Não Feito'
$ws.Cells.Item(113, 3).Value = 'Feito'
$ws.Cells.Item(116, 3).Value = 'Feito'
$ws.Cells.Item(117, 3).Value = 'Feito'
$ws.Cells.Item(121, 3).Value = 'Feito'
$ws.Cells.Item(124, 3).Value = 'Feito'
$ws.Cells.Item(125, 3).Value = 'Feito'
$ws.Cells.Item(127, 3).Value = 'Feito'
$ws.Cells.Item(128, 3).Value = 'Feito'
$ws.Cells.Item(129, 3).Value = 'Feito'
$ws.Cells.Item(130, 3).Value = 'Feito'
$ws.Cells.Item(131, 3).Value = 'Feito'

# --- Fill Mapa / Foto Perfil / Nome no Mapa (columns E, F, G) for rows 106-131 ---
$ws.Cells.Item(106, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(106, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(106, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(107, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(107, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(107, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(108, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(108, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(108, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(109, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(109, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(109, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(110, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(110, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(110, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(111, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(111, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(111, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(112, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(112, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(112, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(113, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(113, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(113, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(114, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(114, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(114, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(115, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(115, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(115, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(116, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(116, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(116, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(117, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(117, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(117, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(118, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/50140/'
$ws.Cells.Item(118, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/50140/file/1149831/blob.-0a54fc10d15569ea9a82a41230fab9bb.png'
$ws.Cells.Item(118, 7).Value = 'Maenuel Severino dos Santos '
$ws.Cells.Item(119, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(119, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(119, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(120, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(120, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(120, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(121, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(121, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(121, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(122, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(122, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(122, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(123, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(123, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(123, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(124, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(124, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(124, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(125, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(125, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(125, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(126, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(126, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(126, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(127, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(127, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(127, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(128, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(128, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(128, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(129, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(129, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(129, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(130, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(130, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(130, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'
$ws.Cells.Item(131, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/38841/'
$ws.Cells.Item(131, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/38841/file/734633/blob-340fa3d2e03e935ef129ffaaaca04b1c.png'
$ws.Cells.Item(131, 7).Value = 'ELTON LEONARDO DE LIMA GALVÃO'

# --- Append new rows 380-409 ---
# Row 380: Artista Elisangela Monteiro / RG
$ws.Cells.Item(380, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(380, 2).Value = 'RG'
$ws.Cells.Item(380, 3).Value = 'Não Feito'
$ws.Cells.Item(380, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(380, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(380, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(380, 7).Value = 'elisangela monteiro de melo costa'

# Row 381: Artista Elisangela Monteiro / CPF
$ws.Cells.Item(381, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(381, 2).Value = 'CPF'
$ws.Cells.Item(381, 3).Value = 'Não Feito'
$ws.Cells.Item(381, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(381, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(381, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(381, 7).Value = 'elisangela monteiro de melo costa'

# Row 382: Artista Elisangela Monteiro / Comprovante de Residência Atual
$ws.Cells.Item(382, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(382, 2).Value = 'Comprovante de Residência Atual'
$ws.Cells.Item(382, 3).Value = 'Não Feito'
$ws.Cells.Item(382, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(382, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(382, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(382, 7).Value = 'elisangela monteiro de melo costa'

# Row 383: Artista Elisangela Monteiro / Comprovações Artísticas
$ws.Cells.Item(383, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(383, 2).Value = 'Comprovações Artísticas'
$ws.Cells.Item(383, 3).Value = 'Não Feito'
$ws.Cells.Item(383, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(383, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(383, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(383, 7).Value = 'elisangela monteiro de melo costa'

# Row 384: Artista Elisangela Monteiro / Currículo Artístico
$ws.Cells.Item(384, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(384, 2).Value = 'Currículo Artístico'
$ws.Cells.Item(384, 3).Value = 'Não Feito'
$ws.Cells.Item(384, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(384, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(384, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(384, 7).Value = 'elisangela monteiro de melo costa'

# Row 385: Artista Elisangela Monteiro / Histórico Atualizado (Ano de Início)
$ws.Cells.Item(385, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(385, 2).Value = 'Histórico Atualizado (Ano de Início)'
$ws.Cells.Item(385, 3).Value = 'Não Feito'
$ws.Cells.Item(385, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(385, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(385, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(385, 7).Value = 'elisangela monteiro de melo costa'

# Row 386: Artista Elisangela Monteiro / Cadastro Mapa Cultural
$ws.Cells.Item(386, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(386, 2).Value = 'Cadastro Mapa Cultural'
$ws.Cells.Item(386, 3).Value = 'Feito'
$ws.Cells.Item(386, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(386, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(386, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(386, 7).Value = 'elisangela monteiro de melo costa'

# Row 387: Artista Elisangela Monteiro / Número Telefone
$ws.Cells.Item(387, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(387, 2).Value = 'Número Telefone'
$ws.Cells.Item(387, 3).Value = 'Não Feito'
$ws.Cells.Item(387, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(387, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(387, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(387, 7).Value = 'elisangela monteiro de melo costa'

# Row 388: Artista Elisangela Monteiro / Cor
$ws.Cells.Item(388, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(388, 2).Value = 'Cor'
$ws.Cells.Item(388, 3).Value = 'Não Feito'
$ws.Cells.Item(388, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(388, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(388, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(388, 7).Value = 'elisangela monteiro de melo costa'

# Row 389: Artista Elisangela Monteiro / Gênero
$ws.Cells.Item(389, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(389, 2).Value = 'Gênero'
$ws.Cells.Item(389, 3).Value = 'Não Feito'
$ws.Cells.Item(389, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(389, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(389, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(389, 7).Value = 'elisangela monteiro de melo costa'

# Row 390: Artista Elisangela Monteiro / Grau de Escolaridade
$ws.Cells.Item(390, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(390, 2).Value = 'Grau de Escolaridade'
$ws.Cells.Item(390, 3).Value = 'Não Feito'
$ws.Cells.Item(390, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(390, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(390, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(390, 7).Value = 'elisangela monteiro de melo costa'

# Row 391: Artista Elisangela Monteiro / Recebe algum benefício do governo?
$ws.Cells.Item(391, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(391, 2).Value = 'Recebe algum benefício do governo?'
$ws.Cells.Item(391, 3).Value = 'Não Feito'
$ws.Cells.Item(391, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(391, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(391, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(391, 7).Value = 'elisangela monteiro de melo costa'

# Row 392: Artista Elisangela Monteiro / Recebeu recursos públicos últimos 5 anos?
$ws.Cells.Item(392, 1).Value = 'Artista Elisangela Monteiro'
$ws.Cells.Item(392, 2).Value = 'Recebeu recursos públicos últimos 5 anos?'
$ws.Cells.Item(392, 3).Value = 'Não Feito'
$ws.Cells.Item(392, 4).Value = '01/09/2023 às 17:55'
$ws.Cells.Item(392, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/16301/'
$ws.Cells.Item(392, 6).Value = 'https://www.mapacultural.pe.gov.br/assets/www.mapacultural.pe.gov.br/img/avatar--agent-2487234669-1693235595.png'
$ws.Cells.Item(392, 7).Value = 'elisangela monteiro de melo costa'

# Row 393: Artista Marcelo Stallone / RG
$ws.Cells.Item(393, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(393, 2).Value = 'RG'
$ws.Cells.Item(393, 3).Value = 'Não Feito'
$ws.Cells.Item(393, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(393, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(393, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(393, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 394: Artista Marcelo Stallone / CPF
$ws.Cells.Item(394, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(394, 2).Value = 'CPF'
$ws.Cells.Item(394, 3).Value = 'Não Feito'
$ws.Cells.Item(394, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(394, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(394, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(394, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 395: Artista Marcelo Stallone / Comprovante de Residência Atual
$ws.Cells.Item(395, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(395, 2).Value = 'Comprovante de Residência Atual'
$ws.Cells.Item(395, 3).Value = 'Não Feito'
$ws.Cells.Item(395, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(395, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(395, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(395, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 396: Artista Marcelo Stallone / Comprovações Artísticas
$ws.Cells.Item(396, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(396, 2).Value = 'Comprovações Artísticas'
$ws.Cells.Item(396, 3).Value = 'Não Feito'
$ws.Cells.Item(396, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(396, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(396, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(396, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 397: Artista Marcelo Stallone / Currículo Artístico
$ws.Cells.Item(397, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(397, 2).Value = 'Currículo Artístico'
$ws.Cells.Item(397, 3).Value = 'Não Feito'
$ws.Cells.Item(397, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(397, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(397, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(397, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 398: Artista Marcelo Stallone / Histórico Atualizado (Ano de Início)
$ws.Cells.Item(398, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(398, 2).Value = 'Histórico Atualizado (Ano de Início)'
$ws.Cells.Item(398, 3).Value = 'Não Feito'
$ws.Cells.Item(398, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(398, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(398, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(398, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 399: Artista Marcelo Stallone / Cadastro Mapa Cultural
$ws.Cells.Item(399, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(399, 2).Value = 'Cadastro Mapa Cultural'
$ws.Cells.Item(399, 3).Value = 'Feito'
$ws.Cells.Item(399, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(399, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(399, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(399, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 400: Artista Marcelo Stallone / Número Telefone
$ws.Cells.Item(400, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(400, 2).Value = 'Número Telefone'
$ws.Cells.Item(400, 3).Value = 'Não Feito'
$ws.Cells.Item(400, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(400, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(400, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(400, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 401: Artista Marcelo Stallone / Cor
$ws.Cells.Item(401, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(401, 2).Value = 'Cor'
$ws.Cells.Item(401, 3).Value = 'Não Feito'
$ws.Cells.Item(401, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(401, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(401, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(401, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 402: Artista Marcelo Stallone / Gênero
$ws.Cells.Item(402, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(402, 2).Value = 'Gênero'
$ws.Cells.Item(402, 3).Value = 'Não Feito'
$ws.Cells.Item(402, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(402, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(402, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(402, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 403: Artista Marcelo Stallone / Grau de Escolaridade
$ws.Cells.Item(403, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(403, 2).Value = 'Grau de Escolaridade'
$ws.Cells.Item(403, 3).Value = 'Não Feito'
$ws.Cells.Item(403, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(403, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(403, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(403, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 404: Artista Marcelo Stallone / Recebe algum benefício do governo?
$ws.Cells.Item(404, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(404, 2).Value = 'Recebe algum benefício do governo?'
$ws.Cells.Item(404, 3).Value = 'Não Feito'
$ws.Cells.Item(404, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(404, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(404, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(404, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 405: Artista Marcelo Stallone / Recebeu recursos públicos últimos 5 anos?
$ws.Cells.Item(405, 1).Value = 'Artista Marcelo Stallone'
$ws.Cells.Item(405, 2).Value = 'Recebeu recursos públicos últimos 5 anos?'
$ws.Cells.Item(405, 3).Value = 'Não Feito'
$ws.Cells.Item(405, 4).Value = '01/09/2023 às 17:57'
$ws.Cells.Item(405, 5).Value = 'https://www.mapacultural.pe.gov.br/agente/39532/'
$ws.Cells.Item(405, 6).Value = 'https://www.mapacultural.pe.gov.br/files/agent/39532/file/751981/blob-d05bf665f1c96a687e15e77668350fab.png'
$ws.Cells.Item(405, 7).Value = 'Marcelo Stallone Monteiro Balbino dos Santos'

# Row 406: AA TESTE / teste 1
$ws.Cells.Item(406, 1).Value = 'AA TESTE'
$ws.Cells.Item(406, 2).Value = 'teste 1'
$ws.Cells.Item(406, 3).Value = 'Não Feito'
$ws.Cells.Item(406, 4).Value = '01/09/2023 às 21:36'
$ws.Cells.Item(406, 5).Value = ''
$ws.Cells.Item(406, 6).Value = ''
$ws.Cells.Item(406, 7).Value = ''

# Row 407: AA TESTE / teste 2
$ws.Cells.Item(407, 1).Value = 'AA TESTE'
$ws.Cells.Item(407, 2).Value = 'teste 2'
$ws.Cells.Item(407, 3).Value = 'Feito'
$ws.Cells.Item(407, 4).Value = '01/09/2023 às 21:36'
$ws.Cells.Item(407, 5).Value = ''
$ws.Cells.Item(407, 6).Value = ''
$ws.Cells.Item(407, 7).Value = ''

# Row 408: AA TESTE / teste 3
$ws.Cells.Item(408, 1).Value = 'AA TESTE'
$ws.Cells.Item(408, 2).Value = 'teste 3'
$ws.Cells.Item(408, 3).Value = 'Não Feito'
$ws.Cells.Item(408, 4).Value = '01/09/2023 às 21:36'
$ws.Cells.Item(408, 5).Value = ''
$ws.Cells.Item(408, 6).Value = ''
$ws.Cells.Item(408, 7).Value = ''

# Row 409: AA TESTE / teste 4
$ws.Cells.Item(409, 1).Value = 'AA TESTE'
$ws.Cells.Item(409, 2).Value = 'teste 4'
$ws.Cells.Item(409, 3).Value = 'Não Feito'
$ws.Cells.Item(409, 4).Value = '01/09/2023 às 21:36'
$ws.Cells.Item(409, 5).Value = ''
$ws.Cells.Item(409, 6).Value = ''
$ws.Cells.Item(409, 7).Value = ''

